# Agregando mas interfaces al proyecto
# Insert a new "Otros" (Others) column right before the existing group
# columns (previously B:D, group1/group2/group3), shifting them to C:E,
# and fill the new column with "---" permissions for every directory row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing group columns (B:D) one slot to the right by inserting
# a fresh column at B. (Formats/values of the old columns travel with
# them automatically.)
$ws.Columns("B:B").Insert()

# --- values first, formats after (PasteSpecial-formats on top of a
#     value already in place keeps style bits like quote-prefix that a
#     plain .Value write on an empty/blank cell would not pick up) ---
$ws.Cells.Item(1, 2).Value = "Otros"
$ws.Cells.Item(2, 2).Value = "---"
$ws.Cells.Item(3, 2).Value = "---"
$ws.Cells.Item(4, 2).Value = "---"
$ws.Cells.Item(5, 2).Value = "---"

# Header formatting: same bold/fill/border/centered look as the other
# header cells (now C1), then add vertical centering on top.
$ws.Range("C1").Copy()
$ws.Cells.Item(1, 2).PasteSpecial(-4122)
$ws.Cells.Item(1, 2).VerticalAlignment = -4108

# Data-cell formatting: match the rest of the permission matrix (now
# C3, a cell that already carries the centered/quote-prefixed style).
$ws.Range("C3").Copy()
$ws.Cells.Item(2, 2).PasteSpecial(-4122)
$ws.Cells.Item(3, 2).PasteSpecial(-4122)
$ws.Cells.Item(4, 2).PasteSpecial(-4122)
$ws.Cells.Item(5, 2).PasteSpecial(-4122)

# Widen the new column. (The shifted group-name column, now C, keeps its
# original best-fit width automatically from the insert - no need to
# touch it.)
$ws.Columns("B:B").ColumnWidth = 18.42578125

# Selection now sits on a single cell instead of the old B1:B5 range.
$ws.Range("B1").Select()
